# Append the latest EUR -> ARS quote as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

# The date text ("2025-09-18") would otherwise be auto-recognized as a date
# value by Excel's normal input parsing, so prefix it with an apostrophe
# (exactly like typing it in the Excel UI) to force it to stay plain text,
# matching every other row in this column. Re-apply the Normal style
# afterwards so the forced-text quote-prefix marker doesn't linger as a
# leftover cell style (the other columns never look numeric, so they don't
# need this treatment).
$ws.Cells.Item($row, 1).Value = "'2025-09-18"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "15:20:43"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,748.7387"
